$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.910.50"
$ws.Range("E2").Value = "  -1.26%  "
$ws.Range("D3").Value = "1.909.46"
$ws.Range("E3").Value = "  -1.14%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.97"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4587"
$ws.Range("E7").Value = "  -0.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3816"
$ws.Range("E8").Value = "  -1.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07727"
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9807"
$ws.Range("E10").Value = "  +0.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.22"
$ws.Range("E11").Value = "  -1.49%  "
$ws.Range("D12").Value = "1.874.85"
$ws.Range("E12").Value = "  -2.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.685"
$ws.Range("E13").Value = "  -1.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.957"
$ws.Range("E14").Value = "  -1.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07064"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "83.85"
$ws.Range("E17").Value = "  -3.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009439"
$ws.Range("E18").Value = "  -2.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.62"
$ws.Range("E19").Value = "  -2.31%  "
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "28.914.37"
$ws.Range("E21").Value = "  -1.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.322"
$ws.Range("E22").Value = "  -2.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.94"
$ws.Range("E23").Value = "  -1.06%  "
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("E25").Value = "  +1.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.03"
$ws.Range("E26").Value = "  -1.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.683"
$ws.Range("E27").Value = "  -1.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "117.54"
$ws.Range("E28").Value = "  -0.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.870"
$ws.Range("E29").Value = "  +1.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09296"
$ws.Range("E30").Value = "  -0.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.8667"
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.102"
$ws.Range("E32").Value = "  -1.27%  "
$ws.Range("E33").Value = "  -3.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.014"
$ws.Range("E34").Value = "  -1.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.05718"
$ws.Range("E35").Value = "  -0.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.152"
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02044"
$ws.Range("E38").Value = "  -1.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.439"
$ws.Range("E39").Value = "  -2.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5498"
$ws.Range("E40").Value = "  -2.67%  "
$ws.Range("E41").Value = "  -1.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.855"
$ws.Range("E42").Value = "  +5.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.351"
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.000002768"
$ws.Range("E44").Value = "  -11.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.162"
$ws.Range("E45").Value = "  +3.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5176"
$ws.Range("E46").Value = "  -1.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.22"
$ws.Range("E47").Value = "  -2.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06905"
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.58"
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.778"
$ws.Range("E50").Value = "  -1.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2879"
$ws.Range("E51").Value = "  -4.03%  "
